$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

$ws.Range("A4").Value = 10
$ws.Range("D4").Value = 4

$ws.Range("A10").Value = 3
$ws.Range("B10").Value = 2
$ws.Range("D10").Value = 4

$ws.Range("A18").Value = 3
$ws.Range("B18").Value = 2
$ws.Range("C18").Value = 10
$ws.Range("F18").Value = 4
